$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.659.09'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '2.998.98'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.04'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.84'
$ws.Range("E6").Value = '  -3.70%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").Value = '2.997.93'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.147'
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.96'
$ws.Range("E11").Value = '  +4.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.29'
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("D16").Value = '3.487.33'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '61.577.64'
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("D19").Value = '2.999.74'
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '454.36'
$ws.Range("E20").Value = '  -3.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.26'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.23'
$ws.Range("E25").Value = '  -5.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.12'
$ws.Range("E26").Value = '  -2.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.51'
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.67'
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("E31").Value = '  -3.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("E32").Value = '  -3.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.59'
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").Value = '0.0₃0821'
$ws.Range("E35").Value = '  +2.35%  '
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("E38").Value = '  -3.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.21'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.40'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("E41").Value = '  +7.51%  '
$ws.Range("E42").Value = '  -3.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '399.44'
$ws.Range("E43").Value = '  -6.33%  '
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.04'
$ws.Range("E45").Value = '  +2.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.268'
$ws.Range("E46").Value = '  -5.10%  '
$ws.Range("D47").Value = '2.721.14'
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.61'
$ws.Range("E48").Value = '  +3.99%  '
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  +0.81%  '
